$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Edipi id and the learner name for row 2
$ws.Range("A2").Value = 1001
$ws.Range("C2").Value = "Nicholas"
$ws.Range("D2").Value = "Fletcher"

# Add an email address in G2 as a mailto hyperlink (Excel auto-creates the
# Hyperlink style/font the first time a hyperlink is inserted)
$ws.Range("G2").Value = "n@gmail.com"
$ws.Range("G2").Hyperlinks.Add($ws.Range("G2"), "mailto:n@gmail.com")

# Restore the active selection to match where the author left off
$ws.Range("G6").Select()
